# AdvisoryBoardTracking.xlsx - "Added email body code for salesmen"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NewSubmissions")
$ws2 = $wb.Worksheets.Item("PendingReSubmissions")

# -- NewSubmissions sheet edits --
# Row 2 (Staten Island Mall job): mark "Check & Letter" as sent
$ws1.Range("F2").Value = "Yes"

# Row 3 (Kings County Hospital job): clear the Job Name / Sub# values and
# record the new email body code ("Boobs") + expected meeting month ("June")
$ws1.Range("C3").Value = ""
$ws1.Range("E3").Value = ""
$ws1.Range("F3").Value = "Boobs"
$ws1.Range("H4").Value = "August"
$ws1.Range("H3").Value = "June"

# Select H3 on NewSubmissions and make it the active sheet/tab
$ws1.Range("H3").Select()
$ws1.Activate()
